$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "vacancies" ---
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "vacancies"

# --- Fill in the new vacancies table (header row + one data row) ---
$ws.Range("A1").Value = "Job Title"
$ws.Range("B1").Value = "Vacancy"
$ws.Range("C1").Value = "Hiring manager"
$ws.Range("D1").Value = "Status"

$ws.Range("A2").Value = "Tester TopGun"
$ws.Range("B2").Value = "Tester01"
# Write D2 before C2 so the shared-string table allocates "All" (14)
# ahead of "Anthony Nolan" (15), matching the source order.
$ws.Range("D2").Value = "All"
$ws.Range("C2").Value = "Anthony Nolan"

# --- Column widths for the new table ---
$ws.Columns.Item(1).ColumnWidth = 13.83
$ws.Columns.Item(2).ColumnWidth = 17.66
$ws.Columns.Item(3).ColumnWidth = 20.16

# --- Make "vacancies" the active sheet/tab and scroll/select like the target view ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C2").Select()
